$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update E426 with new value
$ws.Range("E426").Value = 90

# 2. Copy formatting (A:F) from row 426 down into the new data rows 427:448
#    so the new rows inherit the same per-column cell styles used throughout
#    the data block (A=5, B=13, C=4, D=2, E=2, F=1).
$ws.Range("A426:F426").Copy() | Out-Null
$ws.Range("A427:F448").PasteSpecial(-4122) | Out-Null

# 3. Fill in the new log entries for rows 427-448
$ws.Range("A427").Value = "23.05.2023"
$ws.Range("B427").Value = 0.79166666666666663
$ws.Range("C427").Value = "Practice - Play Track"
$ws.Range("D427").Value = "Code"
$ws.Range("E427").Value = 75
$ws.Range("F427").Value = "Note Track Only Updates On Bar Index Changes (Needs Debugging)"
$ws.Range("A428").Value = "24.05.2023"
$ws.Range("B428").Value = 0.88888888888888884
$ws.Range("C428").Value = "Practice - Track Debugging"
$ws.Range("D428").Value = "Code"
$ws.Range("E428").Value = 40
$ws.Range("F428").Value = "Rewrite How Positions are Calculated, Including Accounting for Note Index Changes (Fixed)"
$ws.Range("A429").Value = "27.05.2023"
$ws.Range("B429").Value = 0.5625
$ws.Range("C429").Value = "Practice - Count Down"
$ws.Range("D429").Value = "Code"
$ws.Range("E429").Value = 75
$ws.Range("F429").Value = "Create a Stylised Cound Down Every Time Play Pressed"
$ws.Range("A430").Value = "27.05.2023"
$ws.Range("B430").Value = 0.65625
$ws.Range("C430").Value = "Practice - Play Audio"
$ws.Range("D430").Value = "Code"
$ws.Range("E430").Value = 45
$ws.Range("F430").Value = "Edit Sample Audio File and Play, Pause and Stop Audio with Track"
$ws.Range("A431").Value = "27.05.2023"
$ws.Range("B431").Value = 0.69097222222222221
$ws.Range("C431").Value = "Practice - Synchronise Audio"
$ws.Range("D431").Value = "Code"
$ws.Range("E431").Value = 20
$ws.Range("F431").Value = "Make Tab Play Together with Audio at the Same Exact Pace"
$ws.Range("A432").Value = "27.05.2023"
$ws.Range("B432").Value = 0.70833333333333337
$ws.Range("C432").Value = "Produce Tabs 5"
$ws.Range("D432").Value = "Other"
$ws.Range("E432").Value = 120
$ws.Range("F432").Value = "Finish Metallica Nothing Else Matters Intro"
$ws.Range("A433").Value = "27.05.2023"
$ws.Range("B433").Value = 0.86458333333333337
$ws.Range("C433").Value = "Practice - Register Strums"
$ws.Range("D433").Value = "Code"
$ws.Range("E433").Value = 65
$ws.Range("F433").Value = "Register Each Strum that the Controller Device Make for Each String"
$ws.Range("A434").Value = "27.05.2023"
$ws.Range("B434").Value = 0.91319444444444453
$ws.Range("C434").Value = "Practice - Accuracy"
$ws.Range("D434").Value = "Code"
$ws.Range("E434").Value = 40
$ws.Range("F434").Value = "Compare Each Current Beat to the Actual Guitar Strums"
$ws.Range("A435").Value = "27.05.2023"
$ws.Range("B435").Value = 0.95833333333333337
$ws.Range("C435").Value = "Software Intro Paragraph"
$ws.Range("D435").Value = "Documentation"
$ws.Range("E435").Value = 30
$ws.Range("F435").Value = "Short Intro for Software Development Paragraph"
$ws.Range("A436").Value = "27.05.2023"
$ws.Range("B436").Value = 0.97916666666666663
$ws.Range("C436").Value = "Materials Paragraph"
$ws.Range("D436").Value = "Documentation"
$ws.Range("E436").Value = 35
$ws.Range("F436").Value = "Write About Audio Cutting and Other Sources"
$ws.Range("A437").Value = "28.05.2023"
$ws.Range("B437").Value = 0
$ws.Range("C437").Value = "Chords and Avatars Paragraph"
$ws.Range("D437").Value = "Documentation"
$ws.Range("E437").Value = 35
$ws.Range("F437").Value = "Document Chord List JSON and Avatar Creation"
$ws.Range("A438").Value = "28.05.2023"
$ws.Range("B438").Value = 2.4305555555555556E-2
$ws.Range("C438").Value = "Digital Design"
$ws.Range("D438").Value = "Documentation"
$ws.Range("E438").Value = 25
$ws.Range("F438").Value = "Some Code Snippets and Explanation plus Snapshot of Neumorphic Elements"
$ws.Range("A439").Value = "28.05.2023"
$ws.Range("B439").Value = 4.1666666666666664E-2
$ws.Range("C439").Value = "Document Registration Page"
$ws.Range("D439").Value = "Documentation"
$ws.Range("E439").Value = 25
$ws.Range("F439").Value = "Snapshot and Explanation"
$ws.Range("A440").Value = "28.05.2023"
$ws.Range("B440").Value = 5.9027777777777783E-2
$ws.Range("C440").Value = "Document Proile Page"
$ws.Range("D440").Value = "Documentation"
$ws.Range("E440").Value = 35
$ws.Range("F440").Value = "Snapshot, Code Snippet from Auto Completion and Explanation"
$ws.Range("A441").Value = "28.05.2023"
$ws.Range("B441").Value = 8.3333333333333329E-2
$ws.Range("C441").Value = "Document Landing Page"
$ws.Range("D441").Value = "Documentation"
$ws.Range("E441").Value = 30
$ws.Range("F441").Value = "Snapshot, Code Snippet about Asyc Await and Explanation"
$ws.Range("A442").Value = "28.05.2023"
$ws.Range("B442").Value = 0.10416666666666667
$ws.Range("C442").Value = "Document Jam Session"
$ws.Range("D442").Value = "Documentation"
$ws.Range("E442").Value = 35
$ws.Range("F442").Value = "Snapshot, Explanation of Guitar Neck Board and Equaliser, Timers etc"
$ws.Range("A443").Value = "28.05.2023"
$ws.Range("B443").Value = 0.1388888888888889
$ws.Range("C443").Value = "Document Chord Page"
$ws.Range("D443").Value = "Documentation"
$ws.Range("E443").Value = 120
$ws.Range("F443").Value = "Filtering, Diagram, Card Generation"
$ws.Range("A444").Value = "28.05.2023"
$ws.Range("B444").Value = 0.22222222222222221
$ws.Range("C444").Value = "Document Compose 1"
$ws.Range("D444").Value = "Documentation"
$ws.Range("E444").Value = 40
$ws.Range("F444").Value = "Basic Intro of Features"
$ws.Range("A445").Value = "29.05.2023"
$ws.Range("B445").Value = 0.5
$ws.Range("C445").Value = "Practice - Accuracy 2"
$ws.Range("D445").Value = "Code"
$ws.Range("E445").Value = 110
$ws.Range("F445").Value = "Read Precision: Every Strum Looks Ahead and Behind for 10 Beats to Find Start Position and Calculates Precision"
$ws.Range("A446").Value = "29.05.2023"
$ws.Range("B446").Value = 0.61805555555555558
$ws.Range("C446").Value = "Practice - Header"
$ws.Range("D446").Value = "Code"
$ws.Range("E446").Value = 90
$ws.Range("F446").Value = "Display Score Title and Score Results"
$ws.Range("A447").Value = "29.05.2023"
$ws.Range("B447").Value = 0.75694444444444453
$ws.Range("C447").Value = "Play - Basics"
$ws.Range("D447").Value = "Code"
$ws.Range("E447").Value = 20
$ws.Range("F447").Value = "Copy All Existing Practice Functionalities to Play and create Play Model"
$ws.Range("A448").Value = "29.05.2023"
$ws.Range("B448").Value = 0.77083333333333337
$ws.Range("C448").Value = "Play - POST, GET All, Get ID"
$ws.Range("D448").Value = "Code"
$ws.Range("E448").Value = 30
$ws.Range("F448").Value = "Save Score After a Play"


# 4. The trailing blank rows (previously 427:428) move down to 449:451 and
#    need the same "blank data row" formatting (A:E only, no F column).
$ws.Range("A427:E428").Copy() | Out-Null
$ws.Range("A449:E451").PasteSpecial(-4122) | Out-Null
$ws.Range("A449:E451").ClearContents() | Out-Null

# 5. The "Total Minutes"/"Total Hours" summary formulas move from rows
#    429:430 down to rows 452:453. Copy their D:E formatting down first.
$ws.Range("D429:E430").Copy() | Out-Null
$ws.Range("D452:E453").PasteSpecial(-4122) | Out-Null

$ws.Range("D452").Value = "Total Minutes"
$ws.Range("E452").Formula = "=SUM(E2:E448)"
$ws.Range("D453").Value = "Total Hours"
$ws.Range("E453").Formula = "=E452 / 60"

# 6. Clear out the old formulas/values that used to live at 429:430 in the
#    D/E columns are already overwritten above via the data fill (step 3),
#    so nothing further is required there.

# 7. Append 23 more blank filler rows (686:708) matching the existing
#    trailing filler-row style (A:B only).
$ws.Range("A685:B685").Copy() | Out-Null
$ws.Range("A686:B708").PasteSpecial(-4122) | Out-Null

# 8. Update the view/selection to match the edited area.
$ws.Range("F448").Select() | Out-Null
